# The "Complexity" column (K) is removed from the dataset sheet. Excel
# shifts every column to its right (product, new-existing, exposed-as-api,
# Interface Weightage) one place to the left, from L:O down to K:N, and
# drops the now-unused "Complexity" shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("K").Delete()

# Mirror Excel's own post-delete selection: the new column K (formerly L)
# becomes the active, fully-selected column.
$ws.Columns("K").Select() | Out-Null
